$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Backlog")
Write-Host $ws1.Range("A1").Style
$ws1.Range("E1").Value = 5
Write-Host $ws1.Range("E1").Style
